# Apply commit "#5: insurance, claim, debt, investment done" to the workbook.
#
# Sheet 6 (具有相當價值之財產): the "otherbonds" property_category label used
# by every antique/jewelry row is renamed to "antique".
#
# Sheet 7 (保險 / insurance): row 1 currently duplicates row-2 data instead of
# holding real column headers, and the data rows are missing the trailing
# category/date/legislator/source/index columns that every other sheet has.
# Row 1 becomes real headers, and rows 2-3 get the full set of columns.
#
# Sheet 8 (事業投資 / investment): same situation as sheet 7 - row 1 becomes
# real headers and row 2 gets the full set of trailing columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 6 - 具有相當價值之財產: otherbonds -> antique (F2:F10)
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
for ($r = 2; $r -le 10; $r++) {
    $ws6.Cells.Item($r, 6).Value = "antique"
}

# ---------------------------------------------------------------------
# Sheet 7 - 保險 (insurance)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item(7)

# Row 1: real headers (company, name, owner, property_category, category,
# date, legislator_name, legislator_id, source_file, index)
$ws7.Cells.Item(1, 2).Value = "company"
$ws7.Cells.Item(1, 3).Value = "name"
$ws7.Cells.Item(1, 4).Value = "owner"
$ws7.Cells.Item(1, 5).Value = "property_category"
$ws7.Cells.Item(1, 6).Value = "category"
$ws7.Cells.Item(1, 7).Value = "date"
$ws7.Cells.Item(1, 8).Value = "legislator_name"
$ws7.Cells.Item(1, 9).Value = "legislator_id"
$ws7.Cells.Item(1, 10).Value = "source_file"
$ws7.Cells.Item(1, 11).Value = "index"

# Row 2: index 101, first insurance policy
$ws7.Cells.Item(2, 2).Value = "南山人壽"
$ws7.Cells.Item(2, 3).Value = "新20年期特別增值分紅養老壽險"
$ws7.Cells.Item(2, 4).Value = "李貴敏"
$ws7.Cells.Item(2, 5).Value = "insurance"
$ws7.Cells.Item(2, 6).Value = "normal"
$ws7.Cells.Item(2, 7).Value = "2012-04-27"
$ws7.Cells.Item(2, 8).Value = "李貴敏"
$ws7.Cells.Item(2, 9).Value = 1739
$ws7.Cells.Item(2, 10).Value = "tmp59331"
$ws7.Cells.Item(2, 11).Value = 101

# Row 3: index 102, second insurance policy
$ws7.Cells.Item(3, 2).Value = "南山人壽"
$ws7.Cells.Item(3, 3).Value = "新20年期終身費特別增值分紅終生保險"
$ws7.Cells.Item(3, 4).Value = "李貴敏"
$ws7.Cells.Item(3, 5).Value = "insurance"
$ws7.Cells.Item(3, 6).Value = "normal"
$ws7.Cells.Item(3, 7).Value = "2012-04-27"
$ws7.Cells.Item(3, 8).Value = "李貴敏"
$ws7.Cells.Item(3, 9).Value = 1739
$ws7.Cells.Item(3, 10).Value = "tmp59331"
$ws7.Cells.Item(3, 11).Value = 102

# Match header / data formatting for the newly added columns to the
# existing B:E columns (bold+border on row 1, plain on rows 2-3).
$ws7.Range("B1").Copy()
$ws7.Range("F1:K1").PasteSpecial(-4122)
$ws7.Range("B2").Copy()
$ws7.Range("F2:K3").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Sheet 8 - 事業投資 (business investment)
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item(8)

# Row 1: real headers (owner, company, address, total, register_date,
# register_reason, property_category, category, date, legislator_name,
# legislator_id, source_file, index)
$ws8.Cells.Item(1, 2).Value = "owner"
$ws8.Cells.Item(1, 3).Value = "company"
$ws8.Cells.Item(1, 4).Value = "address"
$ws8.Cells.Item(1, 5).Value = "total"
$ws8.Cells.Item(1, 6).Value = "register_date"
$ws8.Cells.Item(1, 7).Value = "register_reason"
$ws8.Cells.Item(1, 8).Value = "property_category"
$ws8.Cells.Item(1, 9).Value = "category"
$ws8.Cells.Item(1, 10).Value = "date"
$ws8.Cells.Item(1, 11).Value = "legislator_name"
$ws8.Cells.Item(1, 12).Value = "legislator_id"
$ws8.Cells.Item(1, 13).Value = "source_file"
$ws8.Cells.Item(1, 14).Value = "index"

# Row 2: index 115
$ws8.Cells.Item(2, 2).Value = "李貴敏"
$ws8.Cells.Item(2, 3).Value = "國際通商法律事務所"
$ws8.Cells.Item(2, 4).Value = "臺北市松山區敦化北路168號15樓"
$ws8.Cells.Item(2, 5).Value = 90000000
$ws8.Cells.Item(2, 6).Value = "自民國82年"
$ws8.Cells.Item(2, 7).Value = "合夥"
$ws8.Cells.Item(2, 8).Value = "investment"
$ws8.Cells.Item(2, 9).Value = "normal"
$ws8.Cells.Item(2, 10).Value = "2012-04-27"
$ws8.Cells.Item(2, 11).Value = "李貴敏"
$ws8.Cells.Item(2, 12).Value = 1739
$ws8.Cells.Item(2, 13).Value = "tmp59331"
$ws8.Cells.Item(2, 14).Value = 115

# Match header / data formatting for the newly added columns.
$ws8.Range("B1").Copy()
$ws8.Range("H1:N1").PasteSpecial(-4122)
$ws8.Range("B2").Copy()
$ws8.Range("H2:N2").PasteSpecial(-4122)
